# RF_ASM InstructionList.xlsx - "Documented Utils and Program"
#
# The underlying cell data/text did not change; the sheet was simply
# reopened, column C was widened (losing its old "best fit" auto-size flag
# in favour of an explicit width) and the cursor was left on C5 when the
# file was saved again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Widen the descriptive "Info" column with an explicit width so it no
# longer relies on Excel's auto best-fit (columns A/B keep their
# untouched default width).
$ws.Columns("C:C").ColumnWidth = 72.59244791666667

# Leave the selection where the author left it before saving.
$ws.Range("C5").Select() | Out-Null
